$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the AutoFilter's applied criteria / unhide the filtered-out rows ---
# (the sheet was showing only "Active"/"Armed" origin-state rows; show everything again)
$ws.ShowAllData()

# --- Widen column A slightly to fit the new "Self Recovery" entries ---
$ws.Columns.Item(1).ColumnWidth = 12.14

# --- Append the two new "Self Recovery" rows at the bottom of the table ---
$ws.Range("A30").Value = "Self Recovery"
$ws.Range("A31").Value = "Self Recovery"

# --- Update the view: scroll down a bit and move the selection to A32 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("A32").Select()
